# Re-process the metadata sheet with the newly curated dimensions.
# "municipio-nombre" is promoted from a measure to a dimension (like
# "provincia-nombre" already is), so its concept/type/datatype row
# entries (M2:M4) are updated to mirror column N (provincia-nombre).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = "sdmx-dimension:refArea"
$ws.Range("M3").Value = "dim"
$ws.Range("M4").Value = "URI-Municipio"
